# Insert a new bullet paragraph right after the "...egy-két méret
# megváltoztatása." entry, containing the new "11.19- / 15:20-16:18 – ..."
# log line, and move the trailing _GoBack bookmark onto the new paragraph
# (matching how Word leaves it attached to the last paragraph in the list).

$d = $word.ActiveDocument

# Locate the paragraph to anchor on (by its distinctive trailing text).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*egy-két méret megváltoztatása.*") {
        $anchorIndex = $i
        break
    }
}

$anchor = $d.Paragraphs($anchorIndex)

# Creates a new (empty) paragraph right after the anchor, inheriting the
# anchor paragraph's pPr/rPr (style, numbering level, run formatting).
$null = $anchor.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($anchorIndex + 1)
$target = $newPara.Range

# Build the two runs explicitly via raw OOXML so the text lands as two
# separate <w:r> elements (split right before the time-range text), each
# carrying the same rPr as the anchor paragraph's run.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:szCs w:val="26"/></w:rPr><w:t>11.19-</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">15:20-16:18 – Képek beillesztésének a próbálkozása és swiper módosítások, de a kép magassága mindig túllóg a saját szűlődiv-én. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $target.InsertXML($xml)

# Move the "_GoBack" bookmark off the old (anchor) paragraph and onto the
# end of the freshly inserted paragraph, exactly like the diff shows.
$anchor.Range.Bookmarks("_GoBack").Delete()
$newEnd = $d.Paragraphs($anchorIndex + 1).Range
$newEnd.Collapse(0)
$newEnd.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $newEnd)
